# "Re-ran S(7) w/ slack"
# Update the Y/Z "slack" columns (rows 3-13) on the Strategies sheet with the
# re-run values. Row 14 is a blank spacer row, and rows 16-19 hold
# MIN/MAX/AVERAGE/STDEV formulas over Z3:Z14 that recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strategies")

$newValues = @{
    3  = @{ Y = 0;  Z = -1 }
    4  = @{ Y = 0;  Z = 10 }
    5  = @{ Y = 1;  Z = 5 }
    6  = @{ Y = -1; Z = -24 }
    7  = @{ Y = 0;  Z = -31 }
    8  = @{ Y = 0;  Z = 16 }
    9  = @{ Y = -1; Z = -5 }
    10 = @{ Y = 1;  Z = 15 }
    11 = @{ Y = 0;  Z = 2 }
    12 = @{ Y = -1; Z = -20 }
    13 = @{ Y = -1; Z = -4 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("Y$row").Value = $vals.Y
    $ws.Range("Z$row").Value = $vals.Z
}

# Move the active selection (bottom-right frozen pane) from Y3 to Z14.
$ws.Range("Z14").Select()
